$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed "Price" strings below often look like plain decimal numbers
# (e.g. "1.008", "0.4650", "12.00"). Left alone, Excel would coerce them to
# numeric values on assignment and silently drop significant trailing zeros
# (e.g. "0.4650" -> 0.465). The source data are plain text labels, so each
# target cell is temporarily switched to Text format ("@") before the new
# string is written, then restored to the workbook's default "Normal" style
# so no stray formatting is left behind.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.800.82"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.848.38"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "335.09"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "0.4650"
$ws.Range("E7").Value = "  +1.11%  "

$ws.Range("D8").Value = "0.3867"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "46.78"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").Value = "0.07907"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("D11").Value = "0.9686"
$ws.Range("E11").Value = "  -3.18%  "

$ws.Range("D12").Value = "21.32"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").Value = "1.853.43"
$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").Value = "5.896"
$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").Value = "7.155"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "90.48"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").Value = "0.06615"
$ws.Range("E18").Value = "  -1.27%  "

$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").Value = "17.31"

$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").Value = "27.803.06"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").Value = "5.347"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").Value = "10.84"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").Value = "2.297"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "158.69"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.057.96"
$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("D28").Value = "19.47"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "2.064"
$ws.Range("E29").Value = "  -2.47%  "

$ws.Range("D30").Value = "5.370"
$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("D31").Value = "118.80"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "0.09422"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").Value = "0.9449"
$ws.Range("E33").Value = "  -2.97%  "

$ws.Range("D34").Value = "3.592"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").Value = "5.262"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("D36").Value = "1.327"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").Value = "0.06025"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").Value = "0.02208"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").Value = "8.231"
$ws.Range("E39").Value = "  -0.68%  "

$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").Value = "1.155"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").Value = "0.5809"
$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("D43").Value = "0.1846"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").Value = "10.09"
$ws.Range("E44").Value = "  -2.36%  "

$ws.Range("D45").Value = "1.282"
$ws.Range("E45").Value = "  +3.10%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.00"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5448"
$ws.Range("E47").Value = "  -2.32%  "

$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").Value = "0.06849"
$ws.Range("E49").Value = "  +2.26%  "

$ws.Range("D50").Value = "110.68"
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("D51").Value = "1.007"
$ws.Range("E51").Value = "  -32.37%  "

# Restore the default cell style on the price cells so only the displayed
# value differs from the original workbook (no lingering "@" number format).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "applied crypto price update"
